# Add 2011 and 2012 entries to the main list of yearly MotoGP data files.
$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Column B stores the year as text (same number format as the existing rows).
    $ws.Range("B11:B12").NumberFormat = "@"

    $ws.Range("A11").Value = "C:\Users\zaka\Desktop\MOTOGP\Excels\data\2011.xlsx"
    $ws.Range("B11").Value = "2011"

    $ws.Range("A12").Value = "C:\Users\zaka\Desktop\MOTOGP\Excels\data\2012.xlsx"
    $ws.Range("B12").Value = "2012"

    # Extend the selected range to cover the two new rows, same as a user
    # who re-selected the list after appending data.
    $ws.Activate()
    $ws.Range("A2:B12").Select()
}

# Restore the first sheet as the active/visible tab.
$wb.Worksheets.Item(1).Activate()
